# Add all new reserve people (rows 70-77) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70 - Родионова Александра Владимировна
$ws.Range("A70").Value = "Родионова Александра Владимировна"
$ws.Range("B70").Value = "ИтиАБД"
$ws.Range("C70").Value = "@tg_alessandra_rodionova"
$ws.Range("A70").RowHeight = 61.15

# Row 71 - Бородин Никита Игоревич
$ws.Range("A71").Value = "Бородин Никита Игоревич"
$ws.Range("B71").Value = "ИтиАБД"
$ws.Range("C71").Value = "@call17"
$ws.Range("A71").RowHeight = 37.3

# Row 72 - Логинов Максим Денисович
$ws.Range("A72").Value = "Логинов Максим Денисович"
$ws.Range("B72").Value = "ИтиАБД"
$ws.Range("C72").Value = "@Python_abuser"
$ws.Range("A72").RowHeight = 37.3

# Row 73 - Закревский Константин Сергеевич
$ws.Range("A73").Value = "Закревский Константин Сергеевич"
$ws.Range("B73").Value = "ИтиАБД"
$ws.Range("C73").Value = "@Rigel_125"
$ws.Range("A73").RowHeight = 37.3

# Row 74 - Катлярова Самира Шавкатовна
$ws.Range("A74").Value = "Катлярова Самира Шавкатовна"
$ws.Range("B74").Value = "ИтиАБД"
$ws.Range("C74").Value = "@katmiraa"
$ws.Range("A74").RowHeight = 49.25

# Row 75 - Зарубина Диана Евгеньевна
$ws.Range("A75").Value = "Зарубина Диана Евгеньевна"
$ws.Range("B75").Value = "ИтиАБД"
$ws.Range("C75").Value = "@meowwow66"
$ws.Range("A75").RowHeight = 37.3

# Row 76 - Мамедова Мария Владимировна
$ws.Range("A76").Value = "Мамедова Мария Владимировна"
$ws.Range("B76").Value = "ИтиАБД"
$ws.Range("C76").Value = "@mashvixx1"
$ws.Range("A76").RowHeight = 49.25

# Row 77 - Крылова Мария Григорьевна
$ws.Range("A77").Value = "Крылова Мария Григорьевна"
$ws.Range("B77").Value = "ИтиАБД"
$ws.Range("C77").Value = "@mashkaklyger"
$ws.Range("A77").RowHeight = 49.25

# Update the saved view state (scroll position + active selection)
$ws.Range("D73").Select()
